# Generate Report for Archive
# - Flip the localization status shown on the Overview sheet (zh-cn / de-de
#   columns) and on each per-locale detail sheet's "Status" column from
#   "Ready for handoff" to "In Translation".
# - Narrow the now-shorter "Status" columns so the report reads cleanly.

$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" -------------------

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column widths: shrink the Status columns --------------------------------

$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
